$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Address, $Text)
    $r = $ws.Range($Address)
    $r.NumberFormat = "@"
    $r.Value = $Text
    $r.ClearFormats()
}

# Plain text / link / name / percent-string updates (safe as literal text already)
$ws.Range('D2').Value = '67.785.28'
$ws.Range('E2').Value = '  -1.70%  '
$ws.Range('D3').Value = '3.858.77'
$ws.Range('E3').Value = '  -1.98%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('E5').Value = '  -1.15%  '
$ws.Range('E6').Value = '  +0.81%  '
$ws.Range('D7').Value = '3.858.77'
$ws.Range('E7').Value = '  -1.82%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('E9').Value = '  -0.92%  '
$ws.Range('E10').Value = '  -4.16%  '
$ws.Range('E11').Value = '  -2.00%  '
$ws.Range('E12').Value = '  -2.56%  '
$ws.Range('E13').Value = '  -0.08%  '
$ws.Range('E14').Value = '  -2.43%  '
$ws.Range('D15').Value = '4.505.63'
$ws.Range('E15').Value = '  -1.94%  '
$ws.Range('D16').Value = '3.862.55'
$ws.Range('E16').Value = '  -2.25%  '
$ws.Range('D17').Value = '67.908.26'
$ws.Range('E17').Value = '  -1.61%  '
$ws.Range('E18').Value = '  +3.60%  '
$ws.Range('E19').Value = '  -2.75%  '
$ws.Range('E20').Value = '  -0.51%  '
$ws.Range('E21').Value = '  -1.70%  '
$ws.Range('E22').Value = '  -5.86%  '
$ws.Range('E23').Value = '  +0.52%  '
$ws.Range('E24').Value = '  -5.50%  '
$ws.Range('E25').Value = '  -2.36%  '
$ws.Range('E26').Value = '  -2.45%  '
$ws.Range('E27').Value = '  -1.78%  '
$ws.Range('E28').Value = '  -0.16%  '
$ws.Range('E29').Value = '  -3.12%  '
$ws.Range('E30').Value = '  -1.61%  '
$ws.Range('D31').Value = '4.008.30'
$ws.Range('E31').Value = '  -1.83%  '
$ws.Range('E32').Value = '  -1.98%  '
$ws.Range('E33').Value = '  -4.07%  '
$ws.Range('E34').Value = '  -3.22%  '
$ws.Range('E35').Value = '  -1.88%  '
$ws.Range('D36').Value = '3.823.35'
$ws.Range('E36').Value = '  -1.96%  '
$ws.Range('B37').Value = 'Hedera'
$ws.Range('C37').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('E37').Value = '  -3.14%  '
$ws.Range('B38').Value = 'dogwifhat'
$ws.Range('C38').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('E38').Value = '  +9.73%  '
$ws.Range('E39').Value = '  -2.02%  '
$ws.Range('E40').Value = '  +0.32%  '
$ws.Range('E41').Value = '  -2.21%  '
$ws.Range('E42').Value = '  -0.05%  '
$ws.Range('E43').Value = '  -3.26%  '
$ws.Range('E44').Value = '  +7.51%  '
$ws.Range('E45').Value = '  -3.57%  '
$ws.Range('E46').Value = '  -1.58%  '
$ws.Range('B48').Value = 'OKB'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('E48').Value = '  -2.23%  '
$ws.Range('B49').Value = 'Cosmos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('E49').Value = '  -0.45%  '
$ws.Range('E50').Value = '  +4.76%  '
$ws.Range('E51').Value = '  -0.07%  '

# Numeric-looking price strings: force text storage so formatting/precision (e.g. trailing zeros) is preserved
Set-TextValue 'D5' '597.72'
Set-TextValue 'D6' '169.27'
Set-TextValue 'D9' '0.529'
Set-TextValue 'D11' '6.37'
Set-TextValue 'D12' '0.455'
Set-TextValue 'D13' '0.0000255'
Set-TextValue 'D14' '36.75'
Set-TextValue 'D18' '18.07'
Set-TextValue 'D19' '7.30'
Set-TextValue 'D21' '10.79'
Set-TextValue 'D22' '465.44'
Set-TextValue 'D23' '0.736'
Set-TextValue 'D24' '0.0000159'
Set-TextValue 'D25' '83.01'
Set-TextValue 'D26' '2.22'
Set-TextValue 'D27' '11.99'
Set-TextValue 'D28' '0.999'
Set-TextValue 'D29' '9.94'
Set-TextValue 'D30' '2.94'
Set-TextValue 'D32' '7.66'
Set-TextValue 'D34' '31.06'
Set-TextValue 'D35' '9.39'
Set-TextValue 'D37' '0.104'
Set-TextValue 'D38' '3.66'
Set-TextValue 'D41' '5.88'
Set-TextValue 'D42' '0.999'
Set-TextValue 'D43' '0.312'
Set-TextValue 'D44' '0.000299'
Set-TextValue 'D45' '422.97'
Set-TextValue 'D46' '1.97'
Set-TextValue 'D48' '47.08'
Set-TextValue 'D49' '8.58'
Set-TextValue 'D50' '27.17'
Set-TextValue 'D51' '142.96'
